# UndoRedoSequenceDiagram.pptx update
#
# - Bumps the cached "last saved" date placeholder text (7/29/2017 -> 3/26/2018)
#   on the slide master, every custom (slide) layout, and the notes master.
# - Updates the sequence-diagram labels that describe the old UndoRedoStack /
#   UndoableCommand design to the new UndoRedoCareTaker design:
#     * "x:XYZCommand"  -> ":ModelManager"
#     * "undo()"        -> "resetData(AddressBook)" (and widens/repositions
#       the label box so the longer text still lines up with its arrow)

# ---------------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------------

# PowerPoint's Shape.Left/Top/Width/Height are expressed in points and are
# backed by a single-precision (float32) field internally, so a naive
# EMU/12700.0 conversion can land one EMU short after the round trip through
# that narrower type. Nudge the point value up in tiny increments until the
# float32 round trip reproduces the exact EMU we want.
function EmuToPt($emu) {
    $base = $emu / 12700.0
    for ($i = 0; $i -lt 5000; $i++) {
        $candidate = $base + ($i * 0.000001)
        $asSingle = [single]$candidate
        $backToEmu = [double]$asSingle * 12700.0
        if ([int64]$backToEmu -eq $emu) {
            return $candidate
        }
    }
    return $base
}

# Shapes.Item(<number>) is a positional (1-based) lookup, not a lookup by the
# shape's .Id, and several shapes in this deck share the same .Name, so walk
# the collection and match on .Id explicitly.
function Get-ShapeById($shapes, $id) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $candidate = $shapes.Item($i)
        if ($candidate.Id -eq $id) {
            return $candidate
        }
    }
    return $null
}

# Refresh every "date placeholder" shape in a shape collection whose cached
# text is still the old auto-update date.
function Update-DatePlaceholders($shapes, $oldText, $newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)
        if ($shape.HasTextFrame) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text -eq $oldText) {
                $tr.Text = $newText
            }
        }
    }
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Bump the cached datetimeFigureOut field text everywhere it appears.
# ---------------------------------------------------------------------------

$oldDate = "7/29/2017"
$newDate = "3/26/2018"

Update-DatePlaceholders $p.SlideMaster.Shapes $oldDate $newDate

$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholders $layouts.Item($li).Shapes $oldDate $newDate
}

Update-DatePlaceholders $p.NotesMaster.Shapes $oldDate $newDate

# ---------------------------------------------------------------------------
# 2. Update the sequence-diagram shape text on slide 1.
# ---------------------------------------------------------------------------

$slide = $p.Slides.Item(1)

# "x:XYZCommand" participant header -> ":ModelManager"
$xyzCommandShape = Get-ShapeById $slide.Shapes 84
$xyzCommandShape.TextFrame.TextRange.Text = ":ModelManager"

# "undo()" call label -> "resetData(AddressBook)", resized/repositioned so
# the wider label still sits against its arrow.
$undoShape = Get-ShapeById $slide.Shapes 88
$undoShape.TextFrame.TextRange.Text = "resetData(AddressBook)"
$undoShape.Left = EmuToPt(6142472)
$undoShape.Top = EmuToPt(3416411)
$undoShape.Width = EmuToPt(1615428)
$undoShape.Height = EmuToPt(184666)
